$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column F (shifts old "answer" column to G)
$ws.Columns("F:F").Insert()

# New header for inserted column
$ws.Range("F1").Value = "option5"

# Fill "x" for the new option5 column in each data row (rows 2-13)
$ws.Range("F2:F13").Value = "x"

# Update selection to match target state
$ws.Range("F14").Select()
